$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-11-30T17:06:53"
$ws.Range("U4").Value = 98.62
$ws.Range("V4").Value = 91.94
$ws.Range("W4").Value = 55.34
$ws.Range("X4").Value = 48.76
$ws.Range("Y4").Value = 48.14
$ws.Range("Z4").Value = 52.08
$ws.Range("X5").Value = -41.75
$ws.Range("Y5").Value = -51.81
$ws.Range("U6").Value = -3.94
$ws.Range("V6").Value = -3.77
$ws.Range("W6").Value = -2.27
$ws.Range("X6").Value = -3.26
$ws.Range("Y6").Value = -3.8
$ws.Range("Z6").Value = -1.61
$ws.Range("U9").Value = 300.07
$ws.Range("V9").Value = 300.01
$ws.Range("W9").Value = 134.95
$ws.Range("X9").Value = 51.09
$ws.Range("Y9").Value = 248.39
$ws.Range("Z9").Value = 131.67
$ws.Range("X10").Value = -41.75
$ws.Range("Y10").Value = -51.81
$ws.Range("U11").Value = -3.47
$ws.Range("V11").Value = -2.33
$ws.Range("W11").Value = -0.74
$ws.Range("X11").Value = -0.93
$ws.Range("Y11").Value = -0.62
$ws.Range("Z11").Value = 0.54
$ws.Range("U12").Value = 200.98
$ws.Range("V12").Value = 206.63
$ws.Range("W12").Value = 78.08
$ws.Range("Y12").Value = 197.08
$ws.Range("Z12").Value = 77.43000000000001
$ws.Range("U14").Value = 300.07
$ws.Range("V14").Value = 300.01
$ws.Range("W14").Value = 300.28
$ws.Range("X14").Value = 92.84
$ws.Range("Y14").Value = 300.2
$ws.Range("Z14").Value = 300.4
$ws.Range("W15").Value = 165.32
$ws.Range("Z15").Value = 168.73
$ws.Range("U16").Value = -3.47
$ws.Range("V16").Value = -2.33
$ws.Range("W16").Value = -0.74
$ws.Range("X16").Value = -0.93
$ws.Range("Y16").Value = -0.62
$ws.Range("Z16").Value = 0.54
$ws.Range("U17").Value = 200.98
$ws.Range("V17").Value = 206.63
$ws.Range("W17").Value = 78.08
$ws.Range("Y17").Value = 197.08
$ws.Range("Z17").Value = 77.43000000000001
$ws.Range("U19").Value = 98.05
$ws.Range("V19").Value = 91.5
$ws.Range("W19").Value = 55.13
$ws.Range("X19").Value = 48.41
$ws.Range("Y19").Value = 47.75
$ws.Range("Z19").Value = 51.88
$ws.Range("X20").Value = -41.75
$ws.Range("Y20").Value = -51.81
$ws.Range("U21").Value = -4.51
$ws.Range("V21").Value = -4.21
$ws.Range("W21").Value = -2.48
$ws.Range("X21").Value = -3.61
$ws.Range("Y21").Value = -4.18
$ws.Range("Z21").Value = -1.82
$ws.Range("U24").Value = 98.05
$ws.Range("V24").Value = 91.5
$ws.Range("W24").Value = 55.13
$ws.Range("X24").Value = 48.41
$ws.Range("Y24").Value = 47.75
$ws.Range("Z24").Value = 51.88
$ws.Range("X25").Value = -41.75
$ws.Range("Y25").Value = -51.81
$ws.Range("U26").Value = -4.51
$ws.Range("V26").Value = -4.21
$ws.Range("W26").Value = -2.48
$ws.Range("X26").Value = -3.61
$ws.Range("Y26").Value = -4.18
$ws.Range("Z26").Value = -1.82
$ws.Range("U29").Value = 97.40000000000001
$ws.Range("V29").Value = 90.90000000000001
$ws.Range("W29").Value = 54.76
$ws.Range("X29").Value = 47.89
$ws.Range("Y29").Value = 47.18
$ws.Range("Z29").Value = 51.58
$ws.Range("X30").Value = -41.75
$ws.Range("Y30").Value = -51.81
$ws.Range("U31").Value = -5.16
$ws.Range("V31").Value = -4.82
$ws.Range("W31").Value = -2.85
$ws.Range("X31").Value = -4.12
$ws.Range("Y31").Value = -4.75
$ws.Range("Z31").Value = -2.11
$ws.Range("U34").Value = 339.96
$ws.Range("V34").Value = 345.14
$ws.Range("W34").Value = 353.85
$ws.Range("X34").Value = 95
$ws.Range("Y34").Value = 105.97
$ws.Range("Z34").Value = 557.02
$ws.Range("W35").Value = 165.32
$ws.Range("Z35").Value = 168.73
$ws.Range("U36").Value = -1.71
$ws.Range("V36").Value = -0.38
$ws.Range("W36").Value = 0.88
$ws.Range("X36").Value = 1.24
$ws.Range("Y36").Value = 2.23
$ws.Range("Z36").Value = 2.53
$ws.Range("U37").Value = 239.11
$ws.Range("V37").Value = 249.81
$ws.Range("W37").Value = 130.04
$ws.Range("Z37").Value = 332.07
$ws.Range("U39").Value = 98.62
$ws.Range("V39").Value = 91.94
$ws.Range("W39").Value = 55.34
$ws.Range("X39").Value = 48.76
$ws.Range("Y39").Value = 48.14
$ws.Range("Z39").Value = 52.08
$ws.Range("X40").Value = -41.75
$ws.Range("Y40").Value = -51.81
$ws.Range("U41").Value = -3.94
$ws.Range("V41").Value = -3.77
$ws.Range("W41").Value = -2.27
$ws.Range("X41").Value = -3.26
$ws.Range("Y41").Value = -3.8
$ws.Range("Z41").Value = -1.61
$ws.Range("V44").Value = 96
$ws.Range("W44").Value = 57.67
$ws.Range("X44").Value = 94.05
$ws.Range("Y44").Value = 103.74
$ws.Range("Z44").Value = 54.08
$ws.Range("U46").Value = 0.52
$ws.Range("V46").Value = 0.29
$ws.Range("W46").Value = 0.06
$ws.Range("X46").Value = 0.28
$ws.Range("Z46").Value = 0.38
$ws.Range("U49").Value = 104.87
$ws.Range("V49").Value = 102.37
$ws.Range("W49").Value = 61.35
$ws.Range("X49").Value = 101.04
$ws.Range("Y49").Value = 110.95
$ws.Range("Z49").Value = 57.37
$ws.Range("U51").Value = 2.31
$ws.Range("V51").Value = 6.65
$ws.Range("W51").Value = 3.74
$ws.Range("X51").Value = 7.28
$ws.Range("Y51").Value = 7.21
$ws.Range("Z51").Value = 3.67
$ws.Range("U54").Value = 100.16
$ws.Range("V54").Value = 93.38
$ws.Range("W54").Value = 57.04
$ws.Range("X54").Value = 96.56999999999999
$ws.Range("Y54").Value = 108.18
$ws.Range("Z54").Value = 55.93
$ws.Range("U56").Value = -2.4
$ws.Range("V56").Value = -2.33
$ws.Range("W56").Value = -0.57
$ws.Range("X56").Value = 2.8
$ws.Range("Y56").Value = 4.44
$ws.Range("Z56").Value = 2.24
$ws.Range("U59").Value = 107.28
$ws.Range("V59").Value = 100.01
$ws.Range("W59").Value = 60.07
$ws.Range("X59").Value = 97.78
$ws.Range("Y59").Value = 107.95
$ws.Range("Z59").Value = 55.99
$ws.Range("U61").Value = 4.72
$ws.Range("V61").Value = 4.3
$ws.Range("W61").Value = 2.46
$ws.Range("X61").Value = 4.01
$ws.Range("Y61").Value = 4.21
$ws.Range("Z61").Value = 2.3
$ws.Range("U64").Value = 109.58
$ws.Range("V64").Value = 102.04
$ws.Range("W64").Value = 61.29
$ws.Range("X64").Value = 99.65000000000001
$ws.Range("Y64").Value = 109.9
$ws.Range("Z64").Value = 56.88
$ws.Range("U66").Value = 7.01
$ws.Range("V66").Value = 6.33
$ws.Range("W66").Value = 3.68
$ws.Range("X66").Value = 5.88
$ws.Range("Y66").Value = 6.15
$ws.Range("Z66").Value = 3.19
$ws.Range("V69").Value = 102.92
$ws.Range("W69").Value = 61.81
$ws.Range("X69").Value = 100.61
$ws.Range("Y69").Value = 111.07
$ws.Range("Z69").Value = 57.55
$ws.Range("U71").Value = 7.84
$ws.Range("V71").Value = 7.2
$ws.Range("W71").Value = 4.2
$ws.Range("X71").Value = 6.84
$ws.Range("Y71").Value = 7.33
$ws.Range("Z71").Value = 3.86
$ws.Range("U74").Value = 107.96
$ws.Range("V74").Value = 100.43
$ws.Range("W74").Value = 60.51
$ws.Range("X74").Value = 98.5
$ws.Range("Y74").Value = 108.63
$ws.Range("Z74").Value = 56.23
$ws.Range("V76").Value = 4.72
$ws.Range("W76").Value = 2.9
$ws.Range("X76").Value = 4.73
$ws.Range("Y76").Value = 4.89
$ws.Range("Z76").Value = 2.53
$ws.Range("U79").Value = 108.4
$ws.Range("V79").Value = 100.96
$ws.Range("W79").Value = 60.76
$ws.Range("X79").Value = 98.89
$ws.Range("Y79").Value = 109.05
$ws.Range("Z79").Value = 56.45
$ws.Range("U81").Value = 5.84
$ws.Range("V81").Value = 5.25
$ws.Range("W81").Value = 3.15
$ws.Range("X81").Value = 5.12
$ws.Range("Y81").Value = 5.31
$ws.Range("Z81").Value = 2.76
$ws.Range("U84").Value = 96.76000000000001
$ws.Range("V84").Value = 90.20999999999999
$ws.Range("W84").Value = 55.45
$ws.Range("X84").Value = 95.59
$ws.Range("Y84").Value = 106.95
$ws.Range("Z84").Value = 55.13
$ws.Range("U86").Value = -5.81
$ws.Range("V86").Value = -5.5
$ws.Range("W86").Value = -2.16
$ws.Range("X86").Value = 1.82
$ws.Range("Y86").Value = 3.21
$ws.Range("Z86").Value = 1.43
$ws.Range("U89").Value = 97.40000000000001
$ws.Range("V89").Value = 90.90000000000001
$ws.Range("W89").Value = 54.76
$ws.Range("X89").Value = 47.89
$ws.Range("Y89").Value = 47.18
$ws.Range("Z89").Value = 51.58
$ws.Range("X90").Value = -41.75
$ws.Range("Y90").Value = -51.81
$ws.Range("U91").Value = -5.16
$ws.Range("V91").Value = -4.82
$ws.Range("W91").Value = -2.85
$ws.Range("X91").Value = -4.12
$ws.Range("Y91").Value = -4.75
$ws.Range("Z91").Value = -2.11
